$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above the existing header row, pushing it from row 1 down to row 5,
# to make room for the new "report header" info block.
$ws.Rows("1:4").Insert()

# New label / input cells for the report header block.
$ws.Range("A1").Value = "Sales Agent Name:"
$ws.Range("A2").Value = "From Date:"
$ws.Range("A3").Value = "To Date:"

# Bold the three label cells (new shared cell style).
$ws.Range("A1:A3").Font.Bold = $true

# Give the "From Date" input cell a date number format (new shared cell style)...
$ws.Range("B2").NumberFormat = "mm-dd-yy"

# ...and reuse the exact same style for the "To Date" input cell via a format-only paste,
# instead of independently re-applying the format (which would allocate a duplicate style).
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The longer "Sales Agent Name:" label no longer fits the old column-A width, so resize it.
$ws.Columns("A").AutoFit()

# Restore the cursor/selection to where the author last left it.
$ws.Range("I15").Select()

Write-Host "Done"
